$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header columns, copying the header style from AC1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 70
    $ws.Cells.Item($r, 31).Value = 92
    $ws.Cells.Item($r, 32).Value = 0
}
